# Row permutation: D, J, K, L, M, P columns are reshuffled across rows 2-27
# (weekly price-data shuffle per commit message "Fruta / hortaliza, semanal")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{2=15; 3=18; 4=21; 5=23; 6=3; 7=13; 8=16; 9=26; 10=19; 11=14; 12=10; 13=9; 14=27; 15=4; 16=7; 17=5; 18=17; 19=12; 20=24; 21=11; 22=8; 23=22; 24=20; 25=6; 26=2; 27=25}

$targetCols = @(4, 10, 11, 12, 13, 16)

foreach ($col in $targetCols) {
    $newVals = @{}
    foreach ($destRow in $rowMap.Keys) {
        $srcRow = $rowMap[$destRow]
        $newVals[$destRow] = $ws.Cells.Item($srcRow, $col).Value2
    }
    foreach ($destRow in $newVals.Keys) {
        $ws.Cells.Item($destRow, $col).Value = $newVals[$destRow]
    }
}
